$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# --- Row 13: "Scomm / case study" task gets marked DEFERRED and hidden (filtered out) ---
$ws.Cells.Item(13, 4).Value = "DEFERRED"
$ws.Rows(13).Hidden = $true

# --- A new task is inserted into the visible part of the backlog, and the three tasks
#     that were below it shift down by one slot (within the visible rows), finally
#     revealing one brand-new row (35) at the bottom of the table. The task that used to
#     occupy the inserted slot ("Generic / rework introduction and abstract") is marked
#     DONE and becomes hidden. ---

# Row 26: now the new task
$ws.Cells.Item(26, 1).Value = "Scala.React"
$ws.Cells.Item(26, 2).Value = "Explain the connection between Scala.Reat and CPS transformation"
$ws.Cells.Item(26, 3).Value = 10

# Row 29: now holds what used to be row 26's task, marked DONE & hidden
$ws.Cells.Item(29, 1).Value = "Generic"
$ws.Cells.Item(29, 2).Value = "rework introduction and abstract"
$ws.Cells.Item(29, 4).Value = "DONE"
$ws.Rows(29).Hidden = $true

# Row 32: now holds what used to be row 29's task
$ws.Cells.Item(32, 1).Value = "Scala CPS Plugin"
$ws.Cells.Item(32, 2).Value = "Reference programming with shift/reset"

# Row 34: now holds what used to be row 32's task
$ws.Cells.Item(34, 1).Value = "Scala.Swing"
$ws.Cells.Item(34, 2).Value = "example"
$ws.Cells.Item(34, 3).Value = 15

# Row 35 (new): holds what used to be row 34's task
$ws.Cells.Item(35, 1).Value = "Scomm"
$ws.Cells.Item(35, 2).Value = "spring"
$ws.Cells.Item(35, 3).Value = 30

# Grow the table / autofilter range to include the new row, and refresh the worksheet
# dimension to match.
$tbl.Resize($ws.Range("A1:D35"))

# Restore the active selection to its post-edit location.
$ws.Range("B32").Select() | Out-Null
